# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-21 (serial 45190) to 2023-09-23 (serial 45192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 171 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
